$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, Report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Crime statistics table updates (rows 15-33) ---

# Row 15
$ws.Range("F15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100

# Row 16
$ws.Range("D16").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -50
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = 21.052631578947
$ws.Range("L16").Value = 27.777777777777
$ws.Range("M16").Value = -43.902439024390
$ws.Range("N16").Value = -87.362637362637

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -52.941176470588
$ws.Range("I17").Value = 30
$ws.Range("J17").Value = 43
$ws.Range("K17").Value = -30.232558139534
$ws.Range("L17").Value = -34.782608695652
$ws.Range("M17").Value = -3.225806451612
$ws.Range("N17").Value = -50.819672131147

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = 25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -66.292134831460
$ws.Range("N18").Value = -92.574257425742

# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -30.612244897959
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 133
$ws.Range("K19").Value = -33.082706766917
$ws.Range("L19").Value = -31.538461538461
$ws.Range("M19").Value = 11.25
$ws.Range("N19").Value = -40.666666666666

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 39
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 2.631578947368
$ws.Range("L20").Value = 18.181818181818
$ws.Range("M20").Value = 14.705882352941
$ws.Range("N20").Value = -90.993071593533

# Row 21
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -5.882352941176
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -15.730337078651
$ws.Range("I21").Value = 213
$ws.Range("J21").Value = 266
$ws.Range("K21").Value = -19.924812030075
$ws.Range("L21").Value = -18.390804597701
$ws.Range("M21").Value = -22.826086956521
$ws.Range("N21").Value = -82.780921584478

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 300

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 15.492957746478
$ws.Range("I24").Value = 257
$ws.Range("J24").Value = 216
$ws.Range("K24").Value = 18.981481481481
$ws.Range("L24").Value = 3.212851405622
$ws.Range("M24").Value = 33.160621761658

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 54
$ws.Range("K25").Value = -11.111111111111
$ws.Range("L25").Value = -35.135135135135

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = -18.604651162790
$ws.Range("I26").Value = 87
$ws.Range("J26").Value = 99
$ws.Range("K26").Value = -12.121212121212
$ws.Range("L26").Value = 38.095238095238
$ws.Range("M26").Value = -4.395604395604

# Row 27
$ws.Range("F27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -100

# Row 28
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -47.619047619047
$ws.Range("L28").Value = 57.142857142857

# Row 31
$ws.Range("D31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("J31").Value = 3

# Row 33
$ws.Range("D33").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J33").Value = 1
$ws.Range("J14").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("K33").Value = -100
$ws.Range("K14").Copy()
$ws.Range("K33").PasteSpecial(-4122)

$excel.CutCopyMode = $false
